# edit.ps1 - applies the commit "Add files via upload" to PSdir_PowerPoint.pptx
#
# Summary of the change:
#  1. Slide 2 ("What's to know?"): a new "Objective: ..." paragraph (plus a
#     blank spacer paragraph) is inserted before the existing "1) What's the
#     relationship..." / "2) Can I make any predictions..." paragraphs.
#  2. Two new slides are appended at the end of the deck:
#       - Slide 7 "Observations and Questions for enhancements" (Two Content
#         layout) with explanatory text in the left content placeholder; the
#         right content placeholder is left empty.
#       - Slide 8 "Sources" (Title and Content layout) with a hyperlinked URL
#         pointing at the GSA open-data page.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 2: prepend the "Objective" paragraph (and a blank line) ahead of
#    the existing numbered questions.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$contentShape = $s2.Shapes.Item(2)
$contentRange = $contentShape.TextFrame.TextRange

$objectiveText = "Objective: To help small businesses (government contractors) make informative decisions around competing for bids based off any potential relationships between NAICS and obligation dollars from large prime vendors who are contractually obligated to allocate funds to work with small businesses."

$firstPara = $contentRange.Paragraphs(1, 1)
$firstPara.InsertBefore($objectiveText + "`r" + "`r")

# ---------------------------------------------------------------------------
# 2. New slide 7: "Observations and Questions for enhancements" (Two Content)
# ---------------------------------------------------------------------------
$twoContentLayout = $p.SlideMaster.CustomLayouts.Item(4)
$slide7 = $p.Slides.AddSlide($p.Slides.Count + 1, $twoContentLayout)

$slide7.Shapes.Item(1).TextFrame.TextRange.Text = "Observations and Questions for enhancements "

$leftText = "This clearly didn’t work, but why? A boxplot only shows the range of dollars for each NAICS category, but it’s not enlightening on potential relationships between NAICS and dollars. It just shows with "
$slide7.Shapes.Item(2).TextFrame.TextRange.Text = $leftText

# The third shape (right-hand content placeholder) is intentionally left
# empty, matching the source edit.

# ---------------------------------------------------------------------------
# 3. New slide 8: "Sources" (Title and Content) with a hyperlinked URL
# ---------------------------------------------------------------------------
$titleContentLayout = $p.SlideMaster.CustomLayouts.Item(2)
$slide8 = $p.Slides.AddSlide($p.Slides.Count + 1, $titleContentLayout)

$slide8.Shapes.Item(1).TextFrame.TextRange.Text = "Sources "

$sourceShape = $slide8.Shapes.Item(2)
$sourceRange = $sourceShape.TextFrame.TextRange
$sourceRange.Text = "https://"
$sourceRange.InsertAfter("www.gsa.gov/governmentwide-initiatives/gsa-open-data/gsa-datasets")

$fullUrl = "https://www.gsa.gov/governmentwide-initiatives/gsa-open-data/gsa-datasets"
$sourceFullRange = $sourceShape.TextFrame.TextRange
$action = $sourceFullRange.ActionSettings.Item(1)
$action.Action = 7
$action.Hyperlink.Address = $fullUrl
